# Update quarterly income statement: drop oldest quarter column, shift data left,
# append newest quarter (1402-02-28 / Q4 1401) on the right (column M).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: period ('فصل') headers, columns D (4) .. M (13)
$ws.Cells.Item(8, 4).Value = "فصل سوم منتهی به 1399/09"
$ws.Cells.Item(8, 5).Value = "فصل چهارم منتهی به 1399/12"
$ws.Cells.Item(8, 6).Value = "فصل اول منتهی به 1400/03"
$ws.Cells.Item(8, 7).Value = "فصل دوم منتهی به 1400/06"
$ws.Cells.Item(8, 8).Value = "فصل سوم منتهی به 1400/09"
$ws.Cells.Item(8, 9).Value = "فصل چهارم منتهی به 1400/12"
$ws.Cells.Item(8, 10).Value = "فصل اول منتهی به 1401/03"
$ws.Cells.Item(8, 11).Value = "فصل دوم منتهی به 1401/06"
$ws.Cells.Item(8, 12).Value = "فصل سوم منتهی به 1401/09"
$ws.Cells.Item(8, 13).Value = "فصل چهارم منتهی به 1401/12"

# Row 9: publish dates ('تاریخ انتشار'), columns D (4) .. M (13)
$ws.Cells.Item(9, 4).Value = "1400-10-30 (2)"
$ws.Cells.Item(9, 5).Value = "1401-03-04 (8)"
$ws.Cells.Item(9, 6).Value = "1401-04-30 (2)"
$ws.Cells.Item(9, 7).Value = "1401-08-30 (4)"
$ws.Cells.Item(9, 8).Value = "1401-10-28 (2)"
$ws.Cells.Item(9, 9).Value = "1402-02-28 (7)"
$ws.Cells.Item(9, 10).Value = "1401-04-30"
$ws.Cells.Item(9, 11).Value = "1401-08-30 (2)"
$ws.Cells.Item(9, 12).Value = "1401-10-28"
$ws.Cells.Item(9, 13).Value = "1402-02-28"

# Row 11: فروش (Sales)
$ws.Cells.Item(11, 4).Value = 8359
$ws.Cells.Item(11, 5).Value = 10306
$ws.Cells.Item(11, 6).Value = 12232
$ws.Cells.Item(11, 7).Value = 13466
$ws.Cells.Item(11, 8).Value = 13265
$ws.Cells.Item(11, 9).Value = 15957
$ws.Cells.Item(11, 10).Value = 16022
$ws.Cells.Item(11, 11).Value = 18826
$ws.Cells.Item(11, 12).Value = 19017
$ws.Cells.Item(11, 13).Value = 13150

# Row 12: بهای تمام شده کالای فروش رفته (COGS)
$ws.Cells.Item(12, 4).Value = -3708
$ws.Cells.Item(12, 5).Value = -6035
$ws.Cells.Item(12, 6).Value = -6674
$ws.Cells.Item(12, 7).Value = -6792
$ws.Cells.Item(12, 8).Value = -5838
$ws.Cells.Item(12, 9).Value = -10766
$ws.Cells.Item(12, 10).Value = -10127
$ws.Cells.Item(12, 11).Value = -11225
$ws.Cells.Item(12, 12).Value = -10536
$ws.Cells.Item(12, 13).Value = -7481

# Row 13: سود (زیان) ناخالص (Gross profit)
$ws.Cells.Item(13, 4).Value = 4651
$ws.Cells.Item(13, 5).Value = 4271
$ws.Cells.Item(13, 6).Value = 5558
$ws.Cells.Item(13, 7).Value = 6675
$ws.Cells.Item(13, 8).Value = 7427
$ws.Cells.Item(13, 9).Value = 5191
$ws.Cells.Item(13, 10).Value = 5895
$ws.Cells.Item(13, 11).Value = 7601
$ws.Cells.Item(13, 12).Value = 8481
$ws.Cells.Item(13, 13).Value = 5669

# Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
$ws.Cells.Item(14, 4).Value = -105
$ws.Cells.Item(14, 5).Value = -210
$ws.Cells.Item(14, 6).Value = -204
$ws.Cells.Item(14, 7).Value = -206
$ws.Cells.Item(14, 8).Value = -165
$ws.Cells.Item(14, 9).Value = -570
$ws.Cells.Item(14, 10).Value = -206
$ws.Cells.Item(14, 11).Value = -263
$ws.Cells.Item(14, 12).Value = -204
$ws.Cells.Item(14, 13).Value = -279

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/expense)
$ws.Cells.Item(16, 4).Value = 80
$ws.Cells.Item(16, 5).Value = -390
$ws.Cells.Item(16, 6).Value = 8
$ws.Cells.Item(16, 7).Value = -26
$ws.Cells.Item(16, 8).Value = -3
$ws.Cells.Item(16, 9).Value = -280
$ws.Cells.Item(16, 10).Value = 47
$ws.Cells.Item(16, 11).Value = -51
$ws.Cells.Item(16, 12).Value = 33
$ws.Cells.Item(16, 13).Value = -360

# Row 17: سود (زیان) عملیاتی (Operating profit)
$ws.Cells.Item(17, 4).Value = 4626
$ws.Cells.Item(17, 5).Value = 3672
$ws.Cells.Item(17, 6).Value = 5362
$ws.Cells.Item(17, 7).Value = 6442
$ws.Cells.Item(17, 8).Value = 7259
$ws.Cells.Item(17, 9).Value = 4341
$ws.Cells.Item(17, 10).Value = 5737
$ws.Cells.Item(17, 11).Value = 7287
$ws.Cells.Item(17, 12).Value = 8310
$ws.Cells.Item(17, 13).Value = 5030

# Row 18: هزینه های مالی (Financial expenses)
$ws.Cells.Item(18, 4).Value = -867
$ws.Cells.Item(18, 5).Value = -931
$ws.Cells.Item(18, 6).Value = -1189
$ws.Cells.Item(18, 7).Value = -1366
$ws.Cells.Item(18, 8).Value = -1437
$ws.Cells.Item(18, 9).Value = -1824
$ws.Cells.Item(18, 10).Value = -1625
$ws.Cells.Item(18, 11).Value = -1522
$ws.Cells.Item(18, 12).Value = -1281
$ws.Cells.Item(18, 13).Value = -1194

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-operating income/expense)
$ws.Cells.Item(19, 4).Value = 9
$ws.Cells.Item(19, 5).Value = 26
$ws.Cells.Item(19, 6).Value = 956
$ws.Cells.Item(19, 7).Value = 608
$ws.Cells.Item(19, 8).Value = 13
$ws.Cells.Item(19, 9).Value = 996
$ws.Cells.Item(19, 10).Value = 1582
$ws.Cells.Item(19, 11).Value = 13
$ws.Cells.Item(19, 12).Value = 95
$ws.Cells.Item(19, 13).Value = 739

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit)
$ws.Cells.Item(20, 4).Value = 3769
$ws.Cells.Item(20, 5).Value = 2766
$ws.Cells.Item(20, 6).Value = 5129
$ws.Cells.Item(20, 7).Value = 5684
$ws.Cells.Item(20, 8).Value = 5834
$ws.Cells.Item(20, 9).Value = 3513
$ws.Cells.Item(20, 10).Value = 5694
$ws.Cells.Item(20, 11).Value = 5778
$ws.Cells.Item(20, 12).Value = 7125
$ws.Cells.Item(20, 13).Value = 4574

# Row 21: مالیات (Tax)
$ws.Cells.Item(21, 4).Value = -822
$ws.Cells.Item(21, 5).Value = -710
$ws.Cells.Item(21, 6).Value = -939
$ws.Cells.Item(21, 7).Value = -1139
$ws.Cells.Item(21, 8).Value = -1311
$ws.Cells.Item(21, 9).Value = 293
$ws.Cells.Item(21, 10).Value = -929
$ws.Cells.Item(21, 11).Value = -1298
$ws.Cells.Item(21, 12).Value = -1582
$ws.Cells.Item(21, 13).Value = 716

# Row 22: سود (زیان) خالص عملیات در حال تداوم (Net profit from continuing ops)
$ws.Cells.Item(22, 4).Value = 2947
$ws.Cells.Item(22, 5).Value = 2056
$ws.Cells.Item(22, 6).Value = 4190
$ws.Cells.Item(22, 7).Value = 4545
$ws.Cells.Item(22, 8).Value = 4524
$ws.Cells.Item(22, 9).Value = 3807
$ws.Cells.Item(22, 10).Value = 4766
$ws.Cells.Item(22, 11).Value = 4480
$ws.Cells.Item(22, 12).Value = 5543
$ws.Cells.Item(22, 13).Value = 5289

# Row 24: سود (زیان) خالص (Net profit)
$ws.Cells.Item(24, 4).Value = 2947
$ws.Cells.Item(24, 5).Value = 2056
$ws.Cells.Item(24, 6).Value = 4190
$ws.Cells.Item(24, 7).Value = 4545
$ws.Cells.Item(24, 8).Value = 4524
$ws.Cells.Item(24, 9).Value = 3807
$ws.Cells.Item(24, 10).Value = 4766
$ws.Cells.Item(24, 11).Value = 4480
$ws.Cells.Item(24, 12).Value = 5543
$ws.Cells.Item(24, 13).Value = 5289

# Row 26: سرمایه (Capital)
$ws.Cells.Item(26, 4).Value = 8095
$ws.Cells.Item(26, 5).Value = 9007
$ws.Cells.Item(26, 6).Value = 13202
$ws.Cells.Item(26, 7).Value = 11800
$ws.Cells.Item(26, 8).Value = 10818
$ws.Cells.Item(26, 9).Value = 11182
$ws.Cells.Item(26, 10).Value = 10482
$ws.Cells.Item(26, 11).Value = 14767
$ws.Cells.Item(26, 12).Value = 18939
$ws.Cells.Item(26, 13).Value = 14476

